$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: add effort value in column C (Additional Effort [h])
$ws.Range("C34").Value = 1

# New row 36: date, effort, additional effort blank, task text identical to D35
$ws.Range("A36").Value = 41473
$ws.Range("A36").NumberFormat = 'ddd\ dd/mm/yyyy'
$ws.Range("B36").Value = 1.5
$ws.Range("D36").Value = "Revision of manual"

# Update selection to mirror the author's cursor position after edit
$ws.Range("A37").Select()
